$d = $word.ActiveDocument

# Note: "76÷4=19, 0|75÷4=18, 3" is ordered before "19÷9=2, 1|76÷4=19, 0"
# because the latter produces the text "76÷4=19, 0" which would otherwise be
# re-matched (and incorrectly replaced again) by the former if done later.
$pairs = @(
    @("88÷9=9, 7", "69÷2=34, 1"),
    @("98÷7=14, 0", "22÷7=3, 1"),
    @("34÷9=3, 7", "51÷7=7, 2"),
    @("25÷8=3, 1", "50÷7=7, 1"),
    @("89÷2=44, 1", "85÷9=9, 4"),
    @("47÷2=23, 1", "58÷5=11, 3"),
    @("95÷6=15, 5", "52÷6=8, 4"),
    @("75÷8=9, 3", "46÷6=7, 4"),
    @("56÷3=18, 2", "78÷2=39, 0"),
    @("98÷6=16, 2", "99÷4=24, 3"),
    @("89÷5=17, 4", "34÷6=5, 4"),
    @("42÷6=7, 0", "27÷4=6, 3"),
    @("76÷4=19, 0", "75÷4=18, 3"),
    @("19÷9=2, 1", "76÷4=19, 0"),
    @("72÷9=8, 0", "99÷9=11, 0"),
    @("55÷5=11, 0", "41÷8=5, 1"),
    @("75÷5=15, 0", "48÷3=16, 0"),
    @("77÷8=9, 5", "28÷8=3, 4"),
    @("17÷2=8, 1", "84÷3=28, 0"),
    @("86÷6=14, 2", "24÷8=3, 0"),
    @("57÷5=11, 2", "76÷2=38, 0"),
    @("54÷3=18, 0", "66÷8=8, 2"),
    @("81÷7=11, 4", "34÷2=17, 0"),
    @("38÷2=19, 0", "81÷6=13, 3"),
    @("72÷7=10, 2", "73÷9=8, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
